# library_2651.xlsx edit:
#  - Update s2cDNADate (column A) for rows 24-45 from "01.09.17" to "01.09.18"
#  - Update the sheet's selection/scroll position (view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the s2cDNADate values in column A for rows 24 through 45 ---
# Force the range to remain text (so "01.09.18" is not auto-converted into a
# real date serial number by Excel's automatic data-type detection), then
# restore the cells to the default "Normal" style/General format afterwards
# so the saved styling matches the original (unformatted) cells.
$dateRange = $ws.Range("A24:A45")
$dateRange.NumberFormat = "@"
$dateRange.Value = "01.09.18"
$dateRange.Style = "Normal"

# --- Update the view/selection state ---
$win = $excel.Application.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("A25:A45").Select() | Out-Null
